# 'Wind Gap Simulado' (column E) should be V_sur (column D) - V_norte (column C)
# instead of V_norte - V_sur. This simply flips the sign of every value
# currently stored in column E for rows 3 through 366 (row 2 is 0 and unaffected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 3; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 5).Value = ($dVal - $cVal)
}
